$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student rows added below the header row (A2:A11). Only column A is
# populated, one value per row, mirroring how this data was entered.
$values = @(
    "Creanga",
    "Ion",
    "IS21Z",
    "DA",
    "5.0",
    "Turcanu",
    "Turodr",
    "IS21Z",
    "DA",
    "5.0"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $text = $values[$i]

    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        # Values like "5.0" parse as numbers unless the cell is told to
        # treat input as text first; force that so the literal string,
        # e.g. "5.0", is preserved instead of becoming the number 5.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}
